# Apply updated pl_mw values for Case_4_27 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.155304916976149
$ws.Range("C2").Value = 0.3332793765544579
$ws.Range("D2").Value = 0.6534534883497543
$ws.Range("E2").Value = 0.266713030804496
$ws.Range("G2").Value = 0.002458208114749147
$ws.Range("I2").Value = 0.6592024247168169
$ws.Range("J2").Value = 0.1389749154530335
$ws.Range("N2").Value = 1.126150433234031
$ws.Range("O2").Value = 3.881021786178621
# Row 3
$ws.Range("B3").Value = 1.051205908594852
$ws.Range("C3").Value = 0.3007821263212804
$ws.Range("D3").Value = 0.6430913444166606
$ws.Range("E3").Value = 0.2616146409351572
$ws.Range("G3").Value = 0.002461819266324869
$ws.Range("I3").Value = 0.6643088878418695
$ws.Range("J3").Value = 0.1355858879814917
$ws.Range("N3").Value = 1.137253137991898
$ws.Range("O3").Value = 3.869095063314603
# Row 4
$ws.Range("B4").Value = 0.9874931661412347
$ws.Range("C4").Value = 0.28087848649389
$ws.Range("D4").Value = 0.6370750702869259
$ws.Range("E4").Value = 0.2586319310252136
$ws.Range("G4").Value = 0.002464154742950964
$ws.Range("I4").Value = 0.6678949232929483
$ws.Range("J4").Value = 0.133586503654108
$ws.Range("N4").Value = 1.144556362209549
$ws.Range("O4").Value = 3.864418989603877
# Row 5
$ws.Range("B5").Value = 0.9615821670216746
$ws.Range("C5").Value = 0.2727803032309737
$ws.Range("D5").Value = 0.6347103933201481
$ws.Range("E5").Value = 0.2574535599567582
$ws.Range("G5").Value = 0.002465136289064793
$ws.Range("I5").Value = 0.6694693286050502
$ws.Range("J5").Value = 0.1327921874901037
$ws.Range("N5").Value = 1.147654766529683
$ws.Range("O5").Value = 3.863177587308883
# Row 6
$ws.Range("B6").Value = 0.9572828648594509
$ws.Range("C6").Value = 0.2714363794449071
$ws.Range("D6").Value = 0.6343229944385484
$ws.Range("E6").Value = 0.2572601324387591
$ws.Range("G6").Value = 0.002465301077727984
$ws.Range("I6").Value = 0.6697375790070801
$ws.Range("J6").Value = 0.1326615257431314
$ws.Range("N6").Value = 1.148176641055379
$ws.Range("O6").Value = 3.863011522784035
# Row 7
$ws.Range("B7").Value = 0.9871435075447152
$ws.Range("C7").Value = 0.2807692198405505
$ws.Range("D7").Value = 0.6370428271915785
$ws.Range("E7").Value = 0.2586158889038046
$ws.Range("G7").Value = 0.002464167859532477
$ws.Range("I7").Value = 0.6679156988076862
$ws.Range("J7").Value = 0.1335757084861413
$ws.Range("N7").Value = 1.144597653205928
$ws.Range("O7").Value = 3.864399560520411
# Row 8
$ws.Range("B8").Value = 1.11936973276903
$ws.Range("C8").Value = 0.3220640899215255
$ws.Range("D8").Value = 0.6498087582681364
$ws.Range("E8").Value = 0.2649244206084234
$ws.Range("G8").Value = 0.002459428761075544
$ws.Range("I8").Value = 0.6608694491426341
$ws.Range("J8").Value = 0.1377894378451714
$ws.Range("N8").Value = 1.129877776388206
$ws.Range("O8").Value = 3.87635916886083
# Row 9
$ws.Range("B9").Value = 1.380255089408251
$ws.Range("C9").Value = 0.4034343842626527
$ws.Range("D9").Value = 0.6775927287813204
$ws.Range("E9").Value = 0.2784702942494306
$ws.Range("G9").Value = 0.00245106908422713
$ws.Range("I9").Value = 0.6506391650705581
$ws.Range("J9").Value = 0.146701697608087
$ws.Range("N9").Value = 1.104867537884701
$ws.Range("O9").Value = 3.920888229640639
# Row 10
$ws.Range("B10").Value = 1.572873822318115
$ws.Range("C10").Value = 0.4634565903135126
$ws.Range("D10").Value = 0.699690410205875
$ws.Range("E10").Value = 0.2891439575043435
$ws.Range("G10").Value = 0.00244549039683864
$ws.Range("I10").Value = 0.6453266611246917
$ws.Range("J10").Value = 0.1536496900367581
$ws.Range("N10").Value = 1.088840734150907
$ws.Range("O10").Value = 3.966565105610357
# Row 11
$ws.Range("B11").Value = 1.660703020996777
$ws.Range("C11").Value = 0.4908148870325135
$ws.Range("D11").Value = 0.710111077840935
$ws.Range("E11").Value = 0.2941575942632326
$ws.Range("G11").Value = 0.002443073509957056
$ws.Range("I11").Value = 0.6433920359487075
$ws.Range("J11").Value = 0.1568984150506196
$ws.Range("N11").Value = 1.082059304577513
$ws.Range("O11").Value = 3.990184249498725
# Row 12
$ws.Range("B12").Value = 1.693990534274121
$ws.Range("C12").Value = 0.5011824268323721
$ws.Range("D12").Value = 0.7141101796253508
$ws.Range("E12").Value = 0.2960789394110179
$ws.Range("G12").Value = 0.002442175583298274
$ws.Range("I12").Value = 0.6427290589558865
$ws.Range("J12").Value = 0.1581413482904281
$ws.Range("N12").Value = 1.079564580243229
$ws.Range("O12").Value = 3.999538522907812
# Row 13
$ws.Range("B13").Value = 1.686820222746235
$ws.Range("C13").Value = 0.4989492592060856
$ws.Range("D13").Value = 0.7132465424461998
$ws.Range("E13").Value = 0.2956641289668696
$ws.Range("G13").Value = 0.002442368200003336
$ws.Range("I13").Value = 0.6428687415276784
$ws.Range("J13").Value = 0.157873094219994
$ws.Range("N13").Value = 1.080098605187651
$ws.Range("O13").Value = 3.997505640204963
# Row 14
$ws.Range("B14").Value = 1.663441039889619
$ws.Range("C14").Value = 0.4916676806804503
$ws.Range("D14").Value = 0.7104390233963898
$ws.Range("E14").Value = 0.2943152074535647
$ws.Range("G14").Value = 0.002442999290770742
$ws.Range("I14").Value = 0.6433360953256297
$ws.Range("J14").Value = 0.1570004168670636
$ws.Range("N14").Value = 1.081852593769838
$ws.Range("O14").Value = 3.990945599763847
# Row 15
$ws.Range("B15").Value = 1.649124294140847
$ws.Range("C15").Value = 0.4872084803878352
$ws.Range("D15").Value = 0.7087262435516664
$ws.Range("E15").Value = 0.2934919234796709
$ws.Range("G15").Value = 0.002443388101973735
$ws.Range("I15").Value = 0.6436314391149125
$ws.Range("J15").Value = 0.156467533717958
$ws.Range("N15").Value = 1.082936503415233
$ws.Range("O15").Value = 3.986980860124618
# Row 16
$ws.Range("B16").Value = 1.567138026634893
$ws.Range("C16").Value = 0.4616697283199755
$ws.Range("D16").Value = 0.6990168114578239
$ws.Range("E16").Value = 0.2888194908831352
$ws.Range("G16").Value = 0.00244565077131087
$ws.Range("I16").Value = 0.6454628196506036
$ws.Range("J16").Value = 0.1534391542647739
$ws.Range("N16").Value = 1.089294165990211
$ws.Range("O16").Value = 3.965078846988717
# Row 17
$ws.Range("B17").Value = 1.516894063770678
$ws.Range("C17").Value = 0.4460161890715995
$ws.Range("D17").Value = 0.6931547651768426
$ws.Range("E17").Value = 0.2859936421512401
$ws.Range("G17").Value = 0.002447069745228311
$ws.Range("I17").Value = 0.6467099951117277
$ws.Range("J17").Value = 0.1516039284901041
$ws.Range("N17").Value = 1.093324832959489
$ws.Range("O17").Value = 3.952371437822222
# Row 18
$ws.Range("B18").Value = 1.488014518334182
$ws.Range("C18").Value = 0.4370177720597894
$ws.Range("D18").Value = 0.6898177358400233
$ws.Range("E18").Value = 0.2843831699442134
$ws.Range("G18").Value = 0.002447897284958917
$ws.Range("I18").Value = 0.6474726928887335
$ws.Range("J18").Value = 0.1505566396556617
$ws.Range("N18").Value = 1.095691096378985
$ws.Range("O18").Value = 3.945329728132123
# Row 19
$ws.Range("B19").Value = 1.478239775734892
$ws.Range("C19").Value = 0.4339719422260941
$ws.Range("D19").Value = 0.6886938270276062
$ws.Range("E19").Value = 0.2838404460772708
$ws.Range("G19").Value = 0.002448179433466208
$ws.Range("I19").Value = 0.6477387098589347
$ws.Range("J19").Value = 0.1502034669960608
$ws.Range("N19").Value = 1.096500503886212
$ws.Range("O19").Value = 3.942991374426015
# Row 20
$ws.Range("B20").Value = 1.522240614183772
$ws.Range("C20").Value = 0.4476820102802321
$ws.Range("D20").Value = 0.6937752019473464
$ws.Range("E20").Value = 0.2862929180602194
$ws.Range("G20").Value = 0.002446917515426552
$ws.Range("I20").Value = 0.6465725348360394
$ws.Range("J20").Value = 0.1517984336312708
$ws.Range("N20").Value = 1.092890800782811
$ws.Range("O20").Value = 3.953696490533957
# Row 21
$ws.Range("B21").Value = 1.670307307932092
$ws.Range("C21").Value = 0.4938062535950962
$ws.Range("D21").Value = 0.7112622203543992
$ws.Range("E21").Value = 0.2947107994887048
$ws.Range("G21").Value = 0.002442813455419042
$ws.Range("I21").Value = 0.6431969301916638
$ws.Range("J21").Value = 0.1572563979668189
$ws.Range("N21").Value = 1.081335416413758
$ws.Range("O21").Value = 3.992861295461864
# Row 22
$ws.Range("B22").Value = 1.767243071117377
$ws.Range("C22").Value = 0.5239950080525659
$ws.Range("D22").Value = 0.7230000732983797
$ws.Range("E22").Value = 0.3003452315588007
$ws.Range("G22").Value = 0.002440231989882585
$ws.Range("I22").Value = 0.6413967391075133
$ws.Range("J22").Value = 0.1608976098006991
$ws.Range("N22").Value = 1.074210313162482
$ws.Range("O22").Value = 4.020849636439891
# Row 23
$ws.Range("B23").Value = 1.71549184223835
$ws.Range("C23").Value = 0.5078787516223997
$ws.Range("D23").Value = 0.7167070604847652
$ws.Range("E23").Value = 0.2973258555666192
$ws.Range("G23").Value = 0.002441600574349832
$ws.Range("I23").Value = 0.6423202885946608
$ws.Range("J23").Value = 0.1589474280519454
$ws.Range("N23").Value = 1.07797403291756
$ws.Range("O23").Value = 4.005692292292281
# Row 24
$ws.Range("B24").Value = 1.519823418787951
$ws.Range("C24").Value = 0.4469288893510566
$ws.Range("D24").Value = 0.6934945992658754
$ws.Range("E24").Value = 0.2861575713646545
$ws.Range("G24").Value = 0.002446986301981311
$ws.Range("I24").Value = 0.6466345383527283
$ws.Range("J24").Value = 0.1517104735451369
$ws.Range("N24").Value = 1.093086874155894
$ws.Range("O24").Value = 3.953096612118657
# Row 25
$ws.Range("B25").Value = 1.309510879680374
$ws.Range("C25").Value = 0.3813795952856367
$ws.Range("D25").Value = 0.6697811684486226
$ws.Range("E25").Value = 0.2746794568371129
$ws.Range("G25").Value = 0.00245323126693366
$ws.Range("I25").Value = 0.653020951209065
$ws.Range("J25").Value = 0.1442207448966712
$ws.Range("N25").Value = 1.111220969144533
$ws.Range("O25").Value = 3.906573969084377
